# Auto-generated script to apply scheduled market-data refresh values
# to the Titan_Profits leve-profit calculation sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR).
# Columns H-N hold cached computed values (currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ),
# LeveProfit(NQ/HQ)) that are periodically refreshed from market-board data; no formulas are
# involved, so we just overwrite the cached numbers (or clear cells that no longer apply).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2138
$ws.Range("I40").Value = 2216.6667
$ws.Range("J40").Value = 2085.5557
$ws.Range("K40").Value = 2216.6667
$ws.Range("L40").Value = 2085.5557
$ws.Range("M40").Value = -2041.6667
$ws.Range("N40").Value = -2435.5557
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H137").Value = 90911090
$ws.Range("J137").Value = 2050.75
$ws.Range("L137").Value = 6152.25
$ws.Range("N137").Value = -11252.25
$ws.Range("H138").Value = 5850989
$ws.Range("I138").Value = 1467064
$ws.Range("K138").Value = 4401192
$ws.Range("M138").Value = -4396052
$ws.Range("H141").Value = 2751.875
$ws.Range("I141").Value = 1740
$ws.Range("J141").Value = 4775.625
$ws.Range("K141").Value = 5220
$ws.Range("L141").Value = 14326.875
$ws.Range("M141").Value = -40
$ws.Range("N141").Value = -24686.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2211.7605
$ws.Range("I32").Value = 1487.5574
$ws.Range("J32").Value = 6629.4
$ws.Range("K32").Value = 1487.5574
$ws.Range("L32").Value = 6629.4
$ws.Range("M32").Value = -1200.5574
$ws.Range("N32").Value = -7203.4
$ws.Range("H45").Value = 1155.5555
$ws.Range("I45").Value = 1140
$ws.Range("J45").Value = 1175
$ws.Range("K45").Value = 1140
$ws.Range("L45").Value = 1175
$ws.Range("M45").Value = -763
$ws.Range("N45").Value = -1929
$ws.Range("H61").Value = 2309.2632
$ws.Range("I61").Value = 1633.6666
$ws.Range("J61").Value = 4842.75
$ws.Range("K61").Value = 1633.6666
$ws.Range("L61").Value = 4842.75
$ws.Range("M61").Value = -1421.6666
$ws.Range("N61").Value = -5266.75
$ws.Range("H74").Value = 940.7241
$ws.Range("I74").Value = 953.11536
$ws.Range("K74").Value = 953.11536
$ws.Range("M74").Value = -79.11536000000001
$ws.Range("H77").Value = 940.7241
$ws.Range("I77").Value = 953.11536
$ws.Range("K77").Value = 4765.5768
$ws.Range("M77").Value = -397.5767999999998
$ws.Range("H122").Value = 1676.1111
$ws.Range("I122").Value = 1170
$ws.Range("J122").Value = 1929.1666
$ws.Range("K122").Value = 3510
$ws.Range("L122").Value = 5787.4998
$ws.Range("M122").Value = -1060
$ws.Range("N122").Value = -10687.4998
$ws.Range("H123").Value = 32063
$ws.Range("J123").Value = 32063
$ws.Range("L123").Value = 32063
$ws.Range("N123").Value = -41863
$ws.Range("H136").Value = 2309.2632
$ws.Range("I136").Value = 1633.6666
$ws.Range("J136").Value = 4842.75
$ws.Range("K136").Value = 4900.9998
$ws.Range("L136").Value = 14528.25
$ws.Range("M136").Value = -2350.9998
$ws.Range("N136").Value = -19628.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 3673.6155
$ws.Range("I134").Value = 2099.111
$ws.Range("K134").Value = 6297.333
$ws.Range("M134").Value = -3762.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1525.9656
$ws.Range("I31").Value = 1179.1904
$ws.Range("J31").Value = 2436.25
$ws.Range("K31").Value = 1179.1904
$ws.Range("L31").Value = 2436.25
$ws.Range("M31").Value = -884.1904
$ws.Range("N31").Value = -3026.25
$ws.Range("H34").Value = 1525.9656
$ws.Range("I34").Value = 1179.1904
$ws.Range("J34").Value = 2436.25
$ws.Range("K34").Value = 1179.1904
$ws.Range("L34").Value = 2436.25
$ws.Range("M34").Value = -977.1904
$ws.Range("N34").Value = -2840.25
$ws.Range("H58").Value = 1529.04
$ws.Range("I58").Value = 882.3889
$ws.Range("J58").Value = 3191.8572
$ws.Range("K58").Value = 882.3889
$ws.Range("L58").Value = 3191.8572
$ws.Range("M58").Value = -679.3889
$ws.Range("N58").Value = -3597.8572
$ws.Range("H86").Value = 55557556
$ws.Range("I86").Value = 83334830
$ws.Range("J86").Value = 2999.3333
$ws.Range("K86").Value = 83334830
$ws.Range("L86").Value = 2999.3333
$ws.Range("M86").Value = -83333707
$ws.Range("N86").Value = -5245.3333
$ws.Range("H89").Value = 55557556
$ws.Range("I89").Value = 83334830
$ws.Range("J89").Value = 2999.3333
$ws.Range("K89").Value = 416674150
$ws.Range("L89").Value = 14996.6665
$ws.Range("M89").Value = -416668534
$ws.Range("N89").Value = -26228.6665
$ws.Range("H122").Value = 1412.1177
$ws.Range("I122").Value = 1357.5714
$ws.Range("J122").Value = 1666.6666
$ws.Range("K122").Value = 4072.7142
$ws.Range("L122").Value = 4999.9998
$ws.Range("M122").Value = -1622.7142
$ws.Range("N122").Value = -9899.9998
$ws.Range("H132").Value = 2862.08
$ws.Range("I132").Value = 2322.2856
$ws.Range("J132").Value = 3549.0908
$ws.Range("K132").Value = 6966.8568
$ws.Range("L132").Value = 10647.2724
$ws.Range("M132").Value = -4436.8568
$ws.Range("N132").Value = -15707.2724
$ws.Range("H134").Value = 3514.913
$ws.Range("I134").Value = 1979
$ws.Range("J134").Value = 4502.2856
$ws.Range("K134").Value = 5937
$ws.Range("L134").Value = 13506.8568
$ws.Range("M134").Value = -3402
$ws.Range("N134").Value = -18576.8568
$ws.Range("H136").Value = 1529.04
$ws.Range("I136").Value = 882.3889
$ws.Range("J136").Value = 3191.8572
$ws.Range("K136").Value = 2647.1667
$ws.Range("L136").Value = 9575.571599999999
$ws.Range("M136").Value = -97.16670000000022
$ws.Range("N136").Value = -14675.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 283.7
$ws.Range("J107").Value = 270.7143
$ws.Range("L107").Value = 270.7143
$ws.Range("N107").Value = -4110.7143
$ws.Range("H122").Value = 3705037
$ws.Range("I122").Value = 3705037
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11115111
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -11112661
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 2395.3684
$ws.Range("I126").Value = 2252
$ws.Range("J126").Value = 2461.5386
$ws.Range("K126").Value = 6756
$ws.Range("L126").Value = 7384.6158
$ws.Range("M126").Value = -4286
$ws.Range("N126").Value = -12324.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 458.77274
$ws.Range("I93").Value = 517.9286
$ws.Range("K93").Value = 517.9286
$ws.Range("M93").Value = 730.0714
$ws.Range("H122").Value = 3812.4375
$ws.Range("I122").Value = 3250
$ws.Range("J122").Value = 3999.9167
$ws.Range("K122").Value = 9750
$ws.Range("L122").Value = 11999.7501
$ws.Range("M122").Value = -7300
$ws.Range("N122").Value = -16899.7501
$ws.Range("H132").Value = 3865.0334
$ws.Range("I132").Value = 2983
$ws.Range("J132").Value = 5018.4614
$ws.Range("K132").Value = 8949
$ws.Range("L132").Value = 15055.3842
$ws.Range("M132").Value = -6419
$ws.Range("N132").Value = -20115.3842
$ws.Range("H136").Value = 5907.8696
$ws.Range("I136").Value = 2602
$ws.Range("J136").Value = 9514.272000000001
$ws.Range("K136").Value = 7806
$ws.Range("L136").Value = 28542.816
$ws.Range("M136").Value = -5256
$ws.Range("N136").Value = -33642.81600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H122").Value = 251726
$ws.Range("I122").Value = 502252
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 1506756
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -1504306
$ws.Range("N122").Value = -8500
$ws.Range("H123").Value = 36464
$ws.Range("J123").Value = 36464
$ws.Range("L123").Value = 36464
$ws.Range("N123").Value = -46264
$ws.Range("H132").Value = 33337300
$ws.Range("I132").Value = 100004400
$ws.Range("J132").Value = 3749.8
$ws.Range("K132").Value = 300013200
$ws.Range("L132").Value = 11249.4
$ws.Range("M132").Value = -300010670
$ws.Range("N132").Value = -16309.4
$ws.Range("H136").Value = 19668084
$ws.Range("I136").Value = 20897046
$ws.Range("J136").Value = 4695
$ws.Range("K136").Value = 62691138
$ws.Range("L136").Value = 14085
$ws.Range("M136").Value = -62688588
$ws.Range("N136").Value = -19185

